$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 17.73076433333334
$ws.Cells.Item(2, 8).Value = 53.19229300000001
$ws.Cells.Item(2, 9).Value = 0.004631884691211661
$ws.Cells.Item(2, 10).Value = 0.00463188469121166
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.2457116666666667
$ws.Cells.Item(2, 14).Value = 0.737135
$ws.Cells.Item(2, 15).Value = 0.007550096805344261
$ws.Cells.Item(2, 16).Value = 0.007550096805344259
$ws.Cells.Item(2, 17).Value = 4.356655655617223
$ws.Cells.Item(2, 18).Value = 39.20990090055501
$ws.Cells.Item(2, 19).Value = 0.00003497117780984015
$ws.Cells.Item(2, 20).Value = 0.00003497117780984014

$ws.Cells.Item(3, 7).Value = 17.73076433333334
$ws.Cells.Item(3, 8).Value = 53.19229300000001
$ws.Cells.Item(3, 9).Value = 0.004631884691211661
$ws.Cells.Item(3, 10).Value = 0.00463188469121166
$ws.Cells.Item(3, 14).Value = 94.553567
$ws.Cells.Item(3, 15).Value = 0.9684638283904637
$ws.Cells.Item(3, 16).Value = 0.9684638283904636
$ws.Cells.Item(3, 17).Value = 558.8356711176813
$ws.Cells.Item(3, 18).Value = 5029.521040059131
$ws.Cells.Item(3, 19).Value = 0.004485812780714026
$ws.Cells.Item(3, 20).Value = 0.004485812780714025

$ws.Cells.Item(4, 7).Value = 17.73076433333334
$ws.Cells.Item(4, 8).Value = 53.19229300000001
$ws.Cells.Item(4, 9).Value = 0.004631884691211661
$ws.Cells.Item(4, 10).Value = 0.00463188469121166
$ws.Cells.Item(4, 13).Value = 0.737729
$ws.Cells.Item(4, 14).Value = 2.213187
$ws.Cells.Item(4, 15).Value = 0.02266854253064832
$ws.Cells.Item(4, 16).Value = 0.02266854253064832
$ws.Cells.Item(4, 17).Value = 13.08049904086567
$ws.Cells.Item(4, 18).Value = 117.724491367791
$ws.Cells.Item(4, 19).Value = 0.0001049980751197904
$ws.Cells.Item(4, 20).Value = 0.0001049980751197904

$ws.Cells.Item(5, 7).Value = 17.73076433333334
$ws.Cells.Item(5, 8).Value = 53.19229300000001
$ws.Cells.Item(5, 9).Value = 0.004631884691211661
$ws.Cells.Item(5, 10).Value = 0.00463188469121166
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.042878
$ws.Cells.Item(5, 14).Value = 0.128634
$ws.Cells.Item(5, 15).Value = 0.001317532273543725
$ws.Cells.Item(5, 16).Value = 0.001317532273543725
$ws.Cells.Item(5, 17).Value = 0.7602597130846668
$ws.Cells.Item(5, 18).Value = 6.842337417762001
$ws.Cells.Item(5, 19).Value = 0.000006102657568004474
$ws.Cells.Item(5, 20).Value = 0.000006102657568004472

$ws.Cells.Item(6, 9).Value = 0.9353873458333681
$ws.Cells.Item(6, 10).Value = 0.935387345833368
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.2457116666666667
$ws.Cells.Item(6, 14).Value = 0.737135
$ws.Cells.Item(6, 15).Value = 0.007550096805344261
$ws.Cells.Item(6, 16).Value = 0.007550096805344259
$ws.Cells.Item(6, 17).Value = 879.8061355347988
$ws.Cells.Item(6, 18).Value = 7918.25521981319
$ws.Cells.Item(6, 19).Value = 0.00706226501153596
$ws.Cells.Item(6, 20).Value = 0.007062265011535958

$ws.Cells.Item(7, 9).Value = 0.9353873458333681
$ws.Cells.Item(7, 10).Value = 0.935387345833368
$ws.Cells.Item(7, 14).Value = 94.553567
$ws.Cells.Item(7, 15).Value = 0.9684638283904637
$ws.Cells.Item(7, 16).Value = 0.9684638283904636
$ws.Cells.Item(7, 19).Value = 0.9058888099737784
$ws.Cells.Item(7, 20).Value = 0.9058888099737782

$ws.Cells.Item(8, 9).Value = 0.9353873458333681
$ws.Cells.Item(8, 10).Value = 0.935387345833368
$ws.Cells.Item(8, 13).Value = 0.737729
$ws.Cells.Item(8, 14).Value = 2.213187
$ws.Cells.Item(8, 15).Value = 0.02266854253064832
$ws.Cells.Item(8, 16).Value = 0.02266854253064832
$ws.Cells.Item(8, 17).Value = 2641.545309456008
$ws.Cells.Item(8, 18).Value = 23773.90778510408
$ws.Cells.Item(8, 19).Value = 0.02120386783165395
$ws.Cells.Item(8, 20).Value = 0.02120386783165395

$ws.Cells.Item(9, 9).Value = 0.9353873458333681
$ws.Cells.Item(9, 10).Value = 0.935387345833368
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.042878
$ws.Cells.Item(9, 14).Value = 0.128634
$ws.Cells.Item(9, 15).Value = 0.001317532273543725
$ws.Cells.Item(9, 16).Value = 0.001317532273543725
$ws.Cells.Item(9, 17).Value = 153.5308762145107
$ws.Cells.Item(9, 18).Value = 1381.777885930596
$ws.Cells.Item(9, 19).Value = 0.001232403016399868
$ws.Cells.Item(9, 20).Value = 0.001232403016399868

$ws.Cells.Item(10, 7).Value = 227.2177583333333
$ws.Cells.Item(10, 8).Value = 681.653275
$ws.Cells.Item(10, 9).Value = 0.0593570833501536
$ws.Cells.Item(10, 10).Value = 0.05935708335015359
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.2457116666666667
$ws.Cells.Item(10, 14).Value = 0.737135
$ws.Cells.Item(10, 15).Value = 0.007550096805344261
$ws.Cells.Item(10, 16).Value = 0.007550096805344259
$ws.Cells.Item(10, 17).Value = 55.83005409634722
$ws.Cells.Item(10, 18).Value = 502.470486867125
$ws.Cells.Item(10, 19).Value = 0.0004481517253765477
$ws.Cells.Item(10, 20).Value = 0.0004481517253765476

$ws.Cells.Item(11, 7).Value = 227.2177583333333
$ws.Cells.Item(11, 8).Value = 681.653275
$ws.Cells.Item(11, 9).Value = 0.0593570833501536
$ws.Cells.Item(11, 10).Value = 0.05935708335015359
$ws.Cells.Item(11, 14).Value = 94.553567
$ws.Cells.Item(11, 15).Value = 0.9684638283904637
$ws.Cells.Item(11, 16).Value = 0.9684638283904636
$ws.Cells.Item(11, 17).Value = 7161.416512053547
$ws.Cells.Item(11, 18).Value = 64452.74860848192
$ws.Cells.Item(11, 19).Value = 0.05748518818338161
$ws.Cells.Item(11, 20).Value = 0.05748518818338159

$ws.Cells.Item(12, 7).Value = 227.2177583333333
$ws.Cells.Item(12, 8).Value = 681.653275
$ws.Cells.Item(12, 9).Value = 0.0593570833501536
$ws.Cells.Item(12, 10).Value = 0.05935708335015359
$ws.Cells.Item(12, 13).Value = 0.737729
$ws.Cells.Item(12, 14).Value = 2.213187
$ws.Cells.Item(12, 15).Value = 0.02266854253064832
$ws.Cells.Item(12, 16).Value = 0.02266854253064832
$ws.Cells.Item(12, 17).Value = 167.6251296374917
$ws.Cells.Item(12, 18).Value = 1508.626166737425
$ws.Cells.Item(12, 19).Value = 0.001345538568418194
$ws.Cells.Item(12, 20).Value = 0.001345538568418194

$ws.Cells.Item(13, 7).Value = 227.2177583333333
$ws.Cells.Item(13, 8).Value = 681.653275
$ws.Cells.Item(13, 9).Value = 0.0593570833501536
$ws.Cells.Item(13, 10).Value = 0.05935708335015359
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.042878
$ws.Cells.Item(13, 14).Value = 0.128634
$ws.Cells.Item(13, 15).Value = 0.001317532273543725
$ws.Cells.Item(13, 16).Value = 0.001317532273543725
$ws.Cells.Item(13, 17).Value = 9.742643041816667
$ws.Cells.Item(13, 18).Value = 87.68378737635
$ws.Cells.Item(13, 19).Value = 0.00007820487297725225
$ws.Cells.Item(13, 20).Value = 0.00007820487297725223

$ws.Cells.Item(14, 7).Value = 2.387458333333333
$ws.Cells.Item(14, 8).Value = 7.162374999999999
$ws.Cells.Item(14, 9).Value = 0.0006236861252666267
$ws.Cells.Item(14, 10).Value = 0.0006236861252666266
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.2457116666666667
$ws.Cells.Item(14, 14).Value = 0.737135
$ws.Cells.Item(14, 15).Value = 0.007550096805344261
$ws.Cells.Item(14, 16).Value = 0.007550096805344259
$ws.Cells.Item(14, 17).Value = 0.5866263661805554
$ws.Cells.Item(14, 18).Value = 5.279637295624999
$ws.Cells.Item(14, 19).Value = 0.0000047088906219131
$ws.Cells.Item(14, 20).Value = 0.000004708890621913097

$ws.Cells.Item(15, 7).Value = 2.387458333333333
$ws.Cells.Item(15, 8).Value = 7.162374999999999
$ws.Cells.Item(15, 9).Value = 0.0006236861252666267
$ws.Cells.Item(15, 10).Value = 0.0006236861252666266
$ws.Cells.Item(15, 14).Value = 94.553567
$ws.Cells.Item(15, 15).Value = 0.9684638283904637
$ws.Cells.Item(15, 16).Value = 0.9684638283904636
$ws.Cells.Item(15, 17).Value = 75.24756716018054
$ws.Cells.Item(15, 18).Value = 677.228104441625
$ws.Cells.Item(15, 19).Value = 0.0006040174525897316
$ws.Cells.Item(15, 20).Value = 0.0006040174525897315

$ws.Cells.Item(16, 7).Value = 2.387458333333333
$ws.Cells.Item(16, 8).Value = 7.162374999999999
$ws.Cells.Item(16, 9).Value = 0.0006236861252666267
$ws.Cells.Item(16, 10).Value = 0.0006236861252666266
$ws.Cells.Item(16, 13).Value = 0.737729
$ws.Cells.Item(16, 14).Value = 2.213187
$ws.Cells.Item(16, 15).Value = 0.02266854253064832
$ws.Cells.Item(16, 16).Value = 0.02266854253064832
$ws.Cells.Item(16, 17).Value = 1.761297248791666
$ws.Cells.Item(16, 18).Value = 15.851675239125
$ws.Cells.Item(16, 19).Value = 0.00001413805545638178
$ws.Cells.Item(16, 20).Value = 0.00001413805545638178

$ws.Cells.Item(17, 7).Value = 2.387458333333333
$ws.Cells.Item(17, 8).Value = 7.162374999999999
$ws.Cells.Item(17, 9).Value = 0.0006236861252666267
$ws.Cells.Item(17, 10).Value = 0.0006236861252666266
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.042878
$ws.Cells.Item(17, 14).Value = 0.128634
$ws.Cells.Item(17, 15).Value = 0.001317532273543725
$ws.Cells.Item(17, 16).Value = 0.001317532273543725
$ws.Cells.Item(17, 17).Value = 0.1023694384166667
$ws.Cells.Item(17, 18).Value = 0.9213249457499999
$ws.Cells.Item(17, 19).Value = 0.0000008217265986002152
$ws.Cells.Item(17, 20).Value = 0.0000008217265986002149
